$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1060285
$ws.Range("J17").Value = 1060285
$ws.Range("L17").Value = 3180855
$ws.Range("N17").Value = -3181191

$ws.Range("H40").Value = 2900
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H132").Value = 2099.4119
$ws.Range("I132").Value = 1573.931
$ws.Range("K132").Value = 4721.793
$ws.Range("M132").Value = -2191.793

$ws.Range("H135").Value = 968.82355
$ws.Range("I135").Value = 677.63635
$ws.Range("K135").Value = 6098.72715
$ws.Range("M135").Value = -3563.72715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 36.333332
$ws.Range("I4").Value = 49
$ws.Range("K4").Value = 49
$ws.Range("M4").Value = 67

$ws.Range("H32").Value = 4728.8643
$ws.Range("I32").Value = 2693.16
$ws.Range("J32").Value = 16038.333
$ws.Range("K32").Value = 2693.16
$ws.Range("L32").Value = 16038.333
$ws.Range("M32").Value = -2406.16
$ws.Range("N32").Value = -16612.333

$ws.Range("H45").Value = 2257.4285
$ws.Range("I45").Value = 1943.5
$ws.Range("J45").Value = 2676
$ws.Range("K45").Value = 1943.5
$ws.Range("L45").Value = 2676
$ws.Range("M45").Value = -1566.5
$ws.Range("N45").Value = -3430

$ws.Range("H119").Value = 43725
$ws.Range("J119").Value = 43725
$ws.Range("L119").Value = 43725
$ws.Range("N119").Value = -53401

$ws.Range("H132").Value = 3663.8125
$ws.Range("I132").Value = 3663.8125
$ws.Range("K132").Value = 10991.4375
$ws.Range("M132").Value = -8461.4375

$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1426
$ws.Range("J64").Value = 1234.6666
$ws.Range("L64").Value = 1234.6666
$ws.Range("N64").Value = -1684.6666

$ws.Range("H67").Value = 1426
$ws.Range("J67").Value = 1234.6666
$ws.Range("L67").Value = 1234.6666
$ws.Range("N67").Value = -2794.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 375
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -476

$ws.Range("H3").Value = 10221.333
$ws.Range("I3").Value = 582
$ws.Range("J3").Value = 29500
$ws.Range("K3").Value = 582
$ws.Range("L3").Value = 29500
$ws.Range("M3").Value = -469
$ws.Range("N3").Value = -29726

$ws.Range("H31").Value = 4443.8667
$ws.Range("I31").Value = 2618.4866
$ws.Range("J31").Value = 12886.25
$ws.Range("K31").Value = 2618.4866
$ws.Range("L31").Value = 12886.25
$ws.Range("M31").Value = -2323.4866
$ws.Range("N31").Value = -13476.25

$ws.Range("H34").Value = 4443.8667
$ws.Range("I34").Value = 2618.4866
$ws.Range("J34").Value = 12886.25
$ws.Range("K34").Value = 2618.4866
$ws.Range("L34").Value = 12886.25
$ws.Range("M34").Value = -2416.4866
$ws.Range("N34").Value = -13290.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 10871534
$ws.Range("J7").Value = 20184
$ws.Range("L7").Value = 20184
$ws.Range("N7").Value = -20408

$ws.Range("H8").Value = 10871534
$ws.Range("J8").Value = 20184
$ws.Range("L8").Value = 20184
$ws.Range("N8").Value = -20462

$ws.Range("H13").Value = 784
$ws.Range("I13").Value = 246
$ws.Range("J13").Value = 1214.4
$ws.Range("K13").Value = 246
$ws.Range("L13").Value = 1214.4
$ws.Range("M13").Value = -107
$ws.Range("N13").Value = -1492.4

$ws.Range("H41").Value = 1465
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()

$ws.Range("H122").Value = 5391.6665
$ws.Range("I122").Value = 2937
$ws.Range("K122").Value = 8811
$ws.Range("M122").Value = -6361

$ws.Range("H126").Value = 6093.85
$ws.Range("I126").Value = 6909.875
$ws.Range("K126").Value = 20729.625
$ws.Range("M126").Value = -18259.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3152.2727
$ws.Range("I7").Value = 3265.4
$ws.Range("J7").Value = 3058
$ws.Range("K7").Value = 3265.4
$ws.Range("L7").Value = 3058
$ws.Range("M7").Value = -3153.4
$ws.Range("N7").Value = -3282

$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H30").Value = 2750
$ws.Range("I30").Value = 2750
$ws.Range("K30").Value = 2750
$ws.Range("M30").Value = -2642

$ws.Range("H40").Value = 4254.0625
$ws.Range("I40").Value = 4143
$ws.Range("J40").Value = 4498.4
$ws.Range("K40").Value = 4143
$ws.Range("L40").Value = 4498.4
$ws.Range("M40").Value = -4007
$ws.Range("N40").Value = -4770.4

$ws.Range("H68").Value = 2093.625
$ws.Range("I68").Value = 2392.375
$ws.Range("J68").Value = 1794.875
$ws.Range("K68").Value = 2392.375
$ws.Range("L68").Value = 1794.875
$ws.Range("M68").Value = -1643.375
$ws.Range("N68").Value = -3292.875

$ws.Range("H71").Value = 2093.625
$ws.Range("I71").Value = 2392.375
$ws.Range("J71").Value = 1794.875
$ws.Range("K71").Value = 11961.875
$ws.Range("L71").Value = 8974.375
$ws.Range("M71").Value = -8217.875
$ws.Range("N71").Value = -16462.375

$ws.Range("H126").Value = 3152.2727
$ws.Range("I126").Value = 3265.4
$ws.Range("J126").Value = 3058
$ws.Range("K126").Value = 9796.200000000001
$ws.Range("L126").Value = 9174
$ws.Range("M126").Value = -7326.200000000001
$ws.Range("N126").Value = -14114

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1869.6
$ws.Range("I17").Value = 2087
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 2087
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = -1915
$ws.Range("N17").Value = -1344

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H30").Value = 25001
$ws.Range("J30").Value = 25001
$ws.Range("L30").Value = 25001
$ws.Range("N30").Value = -25215

$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

$ws.Range("H126").Value = 6068.5386
$ws.Range("I126").Value = 6068.5386
$ws.Range("K126").Value = 18205.6158
$ws.Range("M126").Value = -15735.6158
